$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New date header column AM (next day after 31-jul -> 01-ago)
$ws.Range("AM1").NumberFormat = "@"
$ws.Range("AM1").Value = "01-ago"

# New data values for AM2:AM11 (mirrors the pattern of column AL's growth)
$values = @(13, 17, 11, 10, 11, 16, 11, 14, 17, 15)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 39)
    $cell.HorizontalAlignment = -4108
    $cell.NumberFormat = "0"
    $cell.Value = $values[$i]
}

$ws.Range("AP8").Select()
